# Append a new data row (row 48) to the Adafruit IO export sheet, mirroring
# the existing log rows: Timestamp, Feed Key, Value, Latitude, Longitude,
# Elevation.
#
# The new row duplicates the most recent reading (same shape/pattern as the
# immediately preceding row 47): a "temperature" feed reading of "25" with
# "N/A" location fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 48

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"

# "25" must be stored as literal text (matching every other Value/Latitude/
# Longitude/Elevation cell in the sheet, which are all text), not as a
# number. Pre-formatting the cell as Text keeps Excel from auto-coercing the
# numeric-looking string, then resetting the style back to Normal afterwards
# avoids leaving a stray number-format style applied to the cell.
$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "25"
$ws.Range("C$newRow").Style = "Normal"

$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
